$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed B-column values (rows 2-148)
$ws.Cells.Item(2, 2).Value = 61.68613450127155
$ws.Cells.Item(3, 2).Value = 62.42644814919469
$ws.Cells.Item(4, 2).Value = 63.06671941226335
$ws.Cells.Item(5, 2).Value = 63.326829612885
$ws.Cells.Item(6, 2).Value = 63.89707120655552
$ws.Cells.Item(7, 2).Value = 64.64738909296412
$ws.Cells.Item(8, 2).Value = 65.98795705001413
$ws.Cells.Item(9, 2).Value = 66.30809268154846
$ws.Cells.Item(10, 2).Value = 67.70868606951117
$ws.Cells.Item(11, 2).Value = 68.90919468776491
$ws.Cells.Item(12, 2).Value = 70.81
$ws.Cells.Item(14, 2).Value = 70.36
$ws.Cells.Item(15, 2).Value = 71.37
$ws.Cells.Item(16, 2).Value = 72.34999999999999
$ws.Cells.Item(18, 2).Value = 71.73
$ws.Cells.Item(19, 2).Value = 71.53
$ws.Cells.Item(20, 2).Value = 70.95999999999999
$ws.Cells.Item(22, 2).Value = 71.39
$ws.Cells.Item(23, 2).Value = 71.33
$ws.Cells.Item(24, 2).Value = 72.34
$ws.Cells.Item(25, 2).Value = 72.73999999999999
$ws.Cells.Item(26, 2).Value = 73.17
$ws.Cells.Item(27, 2).Value = 73.98
$ws.Cells.Item(28, 2).Value = 73.69
$ws.Cells.Item(29, 2).Value = 74.28
$ws.Cells.Item(30, 2).Value = 74.43000000000001
$ws.Cells.Item(31, 2).Value = 74.45999999999999
$ws.Cells.Item(32, 2).Value = 73.93000000000001
$ws.Cells.Item(34, 2).Value = 75.23
$ws.Cells.Item(35, 2).Value = 75.92
$ws.Cells.Item(36, 2).Value = 75.53
$ws.Cells.Item(38, 2).Value = 76.67
$ws.Cells.Item(39, 2).Value = 77.23999999999999
$ws.Cells.Item(40, 2).Value = 77.98
$ws.Cells.Item(41, 2).Value = 77.62
$ws.Cells.Item(42, 2).Value = 78
$ws.Cells.Item(43, 2).Value = 77.98
$ws.Cells.Item(44, 2).Value = 78.88
$ws.Cells.Item(45, 2).Value = 78.78
$ws.Cells.Item(46, 2).Value = 79.86
$ws.Cells.Item(47, 2).Value = 80.23999999999999
$ws.Cells.Item(49, 2).Value = 82.13
$ws.Cells.Item(50, 2).Value = 82.19
$ws.Cells.Item(51, 2).Value = 81.84
$ws.Cells.Item(52, 2).Value = 83.58
$ws.Cells.Item(53, 2).Value = 83.37
$ws.Cells.Item(54, 2).Value = 83.37
$ws.Cells.Item(55, 2).Value = 83.11
$ws.Cells.Item(57, 2).Value = 83.04000000000001
$ws.Cells.Item(58, 2).Value = 83.61
$ws.Cells.Item(59, 2).Value = 83.38
$ws.Cells.Item(60, 2).Value = 82.25
$ws.Cells.Item(61, 2).Value = 82.37
$ws.Cells.Item(62, 2).Value = 83.12
$ws.Cells.Item(63, 2).Value = 83.22
$ws.Cells.Item(65, 2).Value = 83.58
$ws.Cells.Item(66, 2).Value = 83.31
$ws.Cells.Item(67, 2).Value = 83.28
$ws.Cells.Item(68, 2).Value = 83.41
$ws.Cells.Item(69, 2).Value = 83.91
$ws.Cells.Item(70, 2).Value = 84.56999999999999
$ws.Cells.Item(71, 2).Value = 84.89
$ws.Cells.Item(72, 2).Value = 85.84
$ws.Cells.Item(73, 2).Value = 87.31
$ws.Cells.Item(74, 2).Value = 87.97
$ws.Cells.Item(75, 2).Value = 89.26000000000001
$ws.Cells.Item(76, 2).Value = 89.38
$ws.Cells.Item(77, 2).Value = 89.97
$ws.Cells.Item(78, 2).Value = 90.43000000000001
$ws.Cells.Item(79, 2).Value = 91.09999999999999
$ws.Cells.Item(80, 2).Value = 91.63
$ws.Cells.Item(81, 2).Value = 91.3
$ws.Cells.Item(82, 2).Value = 90.8
$ws.Cells.Item(84, 2).Value = 85.2
$ws.Cells.Item(85, 2).Value = 85.41
$ws.Cells.Item(86, 2).Value = 85.94
$ws.Cells.Item(88, 2).Value = 87.23
$ws.Cells.Item(95, 2).Value = 92.95
$ws.Cells.Item(96, 2).Value = 93.15000000000001
$ws.Cells.Item(98, 2).Value = 93.44
$ws.Cells.Item(99, 2).Value = 93.16
$ws.Cells.Item(100, 2).Value = 92.64
$ws.Cells.Item(101, 2).Value = 93.72
$ws.Cells.Item(102, 2).Value = 94.2
$ws.Cells.Item(103, 2).Value = 94.34
$ws.Cells.Item(104, 2).Value = 95.33
$ws.Cells.Item(105, 2).Value = 95.31
$ws.Cells.Item(106, 2).Value = 95.84
$ws.Cells.Item(107, 2).Value = 96.56999999999999
$ws.Cells.Item(108, 2).Value = 96.36
$ws.Cells.Item(109, 2).Value = 96.86
$ws.Cells.Item(110, 2).Value = 97.38
$ws.Cells.Item(111, 2).Value = 97.88
$ws.Cells.Item(112, 2).Value = 98.73999999999999
$ws.Cells.Item(113, 2).Value = 98.95999999999999
$ws.Cells.Item(114, 2).Value = 99.31
$ws.Cells.Item(115, 2).Value = 99.77
$ws.Cells.Item(116, 2).Value = 101.06
$ws.Cells.Item(117, 2).Value = 101.78
$ws.Cells.Item(118, 2).Value = 102.59
$ws.Cells.Item(119, 2).Value = 103.63
$ws.Cells.Item(120, 2).Value = 103.06
$ws.Cells.Item(121, 2).Value = 103.87
$ws.Cells.Item(122, 2).Value = 103.19
$ws.Cells.Item(123, 2).Value = 103.65
$ws.Cells.Item(124, 2).Value = 104.33
$ws.Cells.Item(125, 2).Value = 104.39
$ws.Cells.Item(126, 2).Value = 104.79
$ws.Cells.Item(127, 2).Value = 104.44
$ws.Cells.Item(128, 2).Value = 102.32
$ws.Cells.Item(129, 2).Value = 93.23999999999999
$ws.Cells.Item(130, 2).Value = 101.33
$ws.Cells.Item(131, 2).Value = 102.31
$ws.Cells.Item(132, 2).Value = 101.68
$ws.Cells.Item(133, 2).Value = 104.07
$ws.Cells.Item(134, 2).Value = 104.16
$ws.Cells.Item(135, 2).Value = 104.72
$ws.Cells.Item(136, 2).Value = 105.43
$ws.Cells.Item(137, 2).Value = 105.59
$ws.Cells.Item(138, 2).Value = 105.9
$ws.Cells.Item(139, 2).Value = 105.53
$ws.Cells.Item(140, 2).Value = 105.03
$ws.Cells.Item(141, 2).Value = 104.95
$ws.Cells.Item(142, 2).Value = 104.95
$ws.Cells.Item(143, 2).Value = 104.66
$ws.Cells.Item(144, 2).Value = 104.55
$ws.Cells.Item(145, 2).Value = 104.28
$ws.Cells.Item(146, 2).Value = 104.3
$ws.Cells.Item(147, 2).Value = 104.49
$ws.Cells.Item(148, 2).Value = 104.81

# Add new row 149 (copy date-cell style from A148, then set values)
$ws.Range("A148").Copy()
$ws.Range("A149").PasteSpecial(-4122)
$ws.Cells.Item(149, 1).Value = 45748
$ws.Cells.Item(149, 2).Value = 104.52
